# Update quarterly sales/cost/profit report: drop oldest quarter (1399/06),
# shift all quarters left, and append the newest quarter (1401/12).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Header rows showing the 10 quarter labels (same rows for every sub-table)
$headerRows = @(8, 31, 54, 69, 92)
$quarterLabels = @("فصل سوم منتهی به 1399/09", "فصل چهارم منتهی به 1399/12", "فصل اول منتهی به 1400/03", "فصل دوم منتهی به 1400/06", "فصل سوم منتهی به 1400/09", "فصل چهارم منتهی به 1400/12", "فصل اول منتهی به 1401/03", "فصل دوم منتهی به 1401/06", "فصل سوم منتهی به 1401/09", "فصل چهارم منتهی به 1401/12")
$quarterCols = @("E", "F", "G", "H", "I", "J", "K", "L", "M", "N")
foreach ($r in $headerRows) {
    for ($i = 0; $i -lt $quarterCols.Length; $i++) {
        $ws.Range($quarterCols[$i] + $r).Value = $quarterLabels[$i]
    }
}

# Data rows: values shift one quarter to the left and a new quarter value is appended
$rowData = @{}
$rowData[11] = @(2740, 2810, 2433, 2885, 3063, "-", 2599, 3279, 3233, 3239)
$rowData[12] = @(1773, 2994, 1919, "-", "-", "-", "-", "-", "-", "-")
$rowData[13] = @(1, 320, 253, 202, 182, "-", 296, 192, 301, 214)
$rowData[14] = @(3754, 4867, 4214, 5071, 4982, "-", 3892, 4661, 4939, 4446)
$rowData[15] = @("-", "-", "-", 474, 750, "-", 365, 516, 568, 458)
$rowData[16] = @(19, 16, 13, 15, 26, "-", 16, 17, 20, 16)
$rowData[17] = @("-", "-", "-", 1556, 2018, "-", 1219, 1618, 1668, 1022)
$rowData[18] = @("-", "-", "-", 1, 52, "-", 95, 103, 162, 142)
$rowData[19] = @("-", "-", "-", 340, 302, "-", 327, 451, 523, 504)
$rowData[20] = @(8287, 11007, 8832, 10544, 11375, 0, 8809, 10837, 11414, 10041)
$rowData[22] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[24] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[25] = @(0, 0, "-", "-", "-", "-", "-", "-", "-", "-")
$rowData[26] = @(0, 0, "-", 0, "-", "-", 0, 0, 0, 0)
$rowData[27] = @(8287, 11007, 8832, 10544, 11375, 0, 8809, 10837, 11414, 10041)
$rowData[34] = @(557307, 586362, 597407, 750421, 890434, "-", 1173609, 1742616, 1828217, 1866027)
$rowData[35] = @(266678, 488854, 329882, "-", "-", "-", "-", "-", "-", "-")
$rowData[36] = @(139, 59040, 55905, 48030, 49321, "-", 156629, 117143, 179178, 123585)
$rowData[37] = @(67569, 87617, 117995, 141992, 139488, "-", 159572, 191130, 202520, 182262)
$rowData[38] = @("-", "-", "-", 79766, 138254, "-", 120141, 203130, 252653, 215441)
$rowData[39] = @(6549, 5652, 5139, 7150, 11918, "-", 11811, 13944, 17166, 13864)
$rowData[40] = @("-", "-", "-", 260930, 407691, "-", 470195, 749402, 869811, 558516)
$rowData[41] = @("-", "-", "-", 71, 10003, "-", 38623, 41881, 66672, 58449)
$rowData[42] = @("-", "-", "-", 118398, 123630, "-", 190150, 295543, 355778, 347981)
$rowData[43] = @(898242, 1227525, 1106328, 1406758, 1770739, 0, 2320730, 3354789, 3771995, 3366125)
$rowData[45] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[47] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[48] = @(0, 0, "-", "-", "-", "-", "-", "-", "-", "-")
$rowData[49] = @(0, 0, "-", 0, "-", "-", 0, 0, 0, 0)
$rowData[50] = @(898242, 1227525, 1106328, 1406758, 1770739, 0, 2320730, 3354789, 3771995, 3366125)
$rowData[57] = @(203396715, 208669751, 245543362, 260111265, 290706497, 299241100, 451561755, 531447393, 565486236, 576112072)
$rowData[58] = @(150410604, 163277889, 171903075, "-", "-", "-", "-", "-", "-", "-")
$rowData[59] = @(139000000, 184500000, 220968379, 237772277, 270994505, 277109453, 529152027, 610119792, 595275748, 577500000)
$rowData[60] = @(17999201, 18002260, 28000712, 28000789, 27998394, 28000356, 41000000, 41006222, 41004252, 40994602)
$rowData[61] = @("-", "-", "-", 168282700, 184338667, 185048465, 329153425, 393662791, 444811620, 470395197)
$rowData[62] = @(344684211, 353250000, 395307692, 476666667, 458384615, 543785714, 738187500, 820235294, 858300000, 866500000)
$rowData[63] = @("-", "-", "-", 167692802, 202027255, 202131397, 385721903, 463165637, 521469424, 546493151)
$rowData[64] = @("-", "-", "-", 71000000, 192365385, 196649123, 406557895, 406611650, 411555556, 411612676)
$rowData[65] = @("-", "-", "-", 348229412, 409370861, 399661871, 581498471, 655305987, 680263862, 690438492)
$rowData[72] = @(-426630, -452719, -438990, -566223, -704266, -717222, -793962, -1098921, -1225033, -1230815)
$rowData[73] = @(-220948, -349716, -270417, "-", "-", "-", "-", "-", "-", "-")
$rowData[74] = @(-546, -46996, -45886, -42208, -46104, -42761, -100063, -53723, -106479, -106148)
$rowData[75] = @(-65364, -72670, -82010, -102882, -114177, -140717, -148909, -173469, -173260, -152082)
$rowData[76] = @("-", "-", "-", -110180, -118507, -100246, -87393, -145187, -160829, -145760)
$rowData[77] = @(-3692, -3866, -3088, -4145, -7033, -8501, -7395, -9273, -12159, -16595)
$rowData[78] = @("-", "-", "-", -378514, -321882, -244550, -333223, -459339, -543679, -295463)
$rowData[79] = @("-", "-", "-", -1144, -5543, -12804, -26041, -32020, -58574, -50483)
$rowData[80] = @("-", "-", "-", -129459, -83125, -91325, -119978, -164918, -214792, -290447)
$rowData[81] = @(-717180, -925967, -840391, -1334755, -1400637, -1358126, -1616964, -2136850, -2494805, -2287793)
$rowData[83] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[85] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[86] = @(0, 0, "-", "-", "-", "-", "-", "-", "-", "-")
$rowData[87] = @(0, 0, "-", 0, "-", 0, 0, 0, 0, 0)
$rowData[88] = @(-717180, -925967, -840391, -1334755, -1400637, -1358126, -1616964, -2136850, -2494805, -2287793)
$rowData[95] = @(130677, 133643, 158417, 184198, 186168, 207433, 379647, 643695, 603184, 635212)
$rowData[96] = @(45730, 139138, 59465, "-", "-", "-", "-", "-", "-", "-")
$rowData[97] = @(-407, 12044, 10019, 5822, 3217, 12938, 56566, 63420, 72699, 17437)
$rowData[98] = @(2205, 14947, 35985, 39110, 25311, 16449, 10663, 17661, 29260, 30180)
$rowData[99] = @("-", "-", "-", 28582, 19747, 14299, 32748, 57943, 91824, 69681)
$rowData[100] = @(2857, 1786, 2051, 3005, 4885, -888, 4416, 4671, 5007, -2731)
$rowData[101] = @("-", "-", "-", 96348, 85809, 47732, 136972, 290063, 326132, 263053)
$rowData[102] = @("-", "-", "-", 313, 4460, -1595, 12582, 9861, 8098, 7966)
$rowData[103] = @("-", "-", "-", 44507, 40505, 19781, 70172, 130625, 140986, 57534)
$rowData[104] = @(181062, 301558, 265937, 401885, 370102, 316149, 703766, 1217939, 1277190, 1078332)
$rowData[106] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[108] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$rowData[109] = @(181062, 301558, 265937, 401885, 370102, 316149, 703766, 1217939, 1277190, 1078332)

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $quarterCols.Length; $i++) {
        $ws.Range($quarterCols[$i] + $r).Value = $vals[$i]
    }
}
